# Applies the cryptos-list refresh described by the commit diff:
# updates the Price (D) and Volume(1h) (E) columns for most rows, and
# swaps the BabyDogeCoin/USDe rows (43/44), including their data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    # Every cell in this sheet is stored as text (coin prices use "."
    # as a thousands separator in some rows), so a purely numeric-looking
    # replacement string (e.g. "568.84") must be forced to stay text
    # instead of being auto-parsed into a Double. A leading apostrophe
    # does that; resetting the style back to Normal afterwards avoids
    # leaving a stray quote-prefix cell style behind.
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

Set-TextValue "D2" "62.420.77"
Set-TextValue "E2" "  -1.46%  "
Set-TextValue "D3" "2.446.66"
Set-TextValue "E3" "  -1.18%  "
Set-TextValue "D5" "568.84"
Set-TextValue "E5" "  -1.63%  "
Set-TextValue "D6" "145.78"
Set-TextValue "E6" "  -0.89%  "
Set-TextValue "E7" "  +0.10%  "
Set-TextValue "D8" "0.530"
Set-TextValue "E8" "  -2.29%  "
Set-TextValue "D9" "0.110"
Set-TextValue "E9" "  -1.43%  "
Set-TextValue "E10" "  +0.03%  "
Set-TextValue "D11" "5.23"
Set-TextValue "E11" "  -1.22%  "
Set-TextValue "E12" "  -2.00%  "
Set-TextValue "D13" "28.58"
Set-TextValue "E13" "  -1.97%  "
Set-TextValue "D14" "0.0000174"
Set-TextValue "E14" "  -3.76%  "
Set-TextValue "D15" "2.897.28"
Set-TextValue "E15" "  -0.73%  "
Set-TextValue "D16" "62.428.67"
Set-TextValue "E16" "  -1.22%  "
Set-TextValue "D17" "2.450.91"
Set-TextValue "E17" "  -0.69%  "
Set-TextValue "D18" "7.85"
Set-TextValue "E18" "  -1.10%  "
Set-TextValue "D19" "10.77"
Set-TextValue "E19" "  -3.23%  "
Set-TextValue "D20" "322.54"
Set-TextValue "E20" "  -2.65%  "
Set-TextValue "E21" "  -0.15%  "
Set-TextValue "E22" "  -4.21%  "
Set-TextValue "E23" "  -0.13%  "
Set-TextValue "D24" "9.97"
Set-TextValue "E24" "  +8.50%  "
Set-TextValue "D25" "65.13"
Set-TextValue "E25" "  -2.16%  "
Set-TextValue "D26" "634.84"
Set-TextValue "E26" "  -6.12%  "
Set-TextValue "E27" "  -0.35%  "
Set-TextValue "D28" "0.0₃0958"
Set-TextValue "E28" "  -5.97%  "
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  -4.15%  "
Set-TextValue "E30" "  -3.87%  "
Set-TextValue "D31" "7.86"
Set-TextValue "E31" "  -4.22%  "
Set-TextValue "E32" "  -4.24%  "
Set-TextValue "E33" "  -5.76%  "
Set-TextValue "D34" "0.998"
Set-TextValue "E34" "  -0.01%  "
Set-TextValue "D35" "1.49"
Set-TextValue "E35" "  -4.67%  "
Set-TextValue "E36" "  -2.35%  "
Set-TextValue "D37" "151.12"
Set-TextValue "E37" "  -1.63%  "
Set-TextValue "D38" "0.365"
Set-TextValue "E38" "  -2.81%  "
Set-TextValue "D39" "18.48"
Set-TextValue "E39" "  -2.26%  "
Set-TextValue "D40" "5.25"
Set-TextValue "E40" "  -6.11%  "
Set-TextValue "E41" "  -4.30%  "
Set-TextValue "E42" "  -4.27%  "
Set-TextValue "B43" "USDe"
Set-TextValue "C43" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D43" "0.999"
Set-TextValue "E43" "  -0.08%  "
Set-TextValue "B44" "BabyDogeCoin"
Set-TextValue "C44" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D44" "0.0₆0306"
Set-TextValue "E44" "  +1.73%  "
Set-TextValue "D45" "151.98"
Set-TextValue "E45" "  +3.31%  "
Set-TextValue "E46" "  +1.02%  "
Set-TextValue "D47" "3.53"
Set-TextValue "E47" "  -3.12%  "
Set-TextValue "D48" "0.602"
Set-TextValue "E48" "  -1.15%  "
Set-TextValue "D49" "20.01"
Set-TextValue "E49" "  -4.25%  "
Set-TextValue "D50" "0.0502"
Set-TextValue "E50" "  -3.21%  "
Set-TextValue "D51" "0.0902"
Set-TextValue "E51" "  -2.39%  "
